# Apply the workbook restructuring described by the diff:
#  - keep the original "分数表" sheet (long/tidy layout) with its data unchanged
#  - add a new first sheet "xlwings转换后" with the data pivoted into a wide layout
#  - add a new last sheet "转换后" with the same pivoted layout
#  - final sheet order: xlwings转换后, 分数表, 转换后

$wb = $excel.ActiveWorkbook

# --- locate / normalize the original sheet -------------------------------
$orig = $wb.Worksheets.Item(1)
$orig.Name = "分数表"

# --- create the trailing "转换后" sheet (goes right after 分数表) ---------
$afterRef = $wb.Worksheets.Item("分数表")
$newAfter = $wb.Worksheets.Add($null, $afterRef)
$newAfter.Name = "转换后"

# --- create the leading "xlwings转换后" sheet (goes right before 分数表) --
$beforeRef = $wb.Worksheets.Item("分数表")
$newBefore = $wb.Worksheets.Add($beforeRef)
$newBefore.Name = "xlwings转换后"

# Sheet insertion shifts the internal index of earlier handles, so re-fetch
# fresh, stable references by name now that the sheet set/order is final.
$orig = $wb.Worksheets.Item("分数表")
$converted = $wb.Worksheets.Item("转换后")
$xlwings = $wb.Worksheets.Item("xlwings转换后")

# --- pivoted (wide) data shared by the two new sheets ---------------------
# columns: 姓名 | 语文 | 数学 | 英语
$wide = @(
    @("姓名", "语文", "数学", "英语"),
    @("小张", 96, 100, 97),
    @("小曾", 96, $null, $null),
    @("小江", 98, 99, 100),
    @("李飞", 95, $null, 98)
)

function Fill-Wide($sheet) {
    for ($r = 0; $r -lt $wide.Length; $r++) {
        $row = $wide[$r]
        for ($c = 0; $c -lt $row.Length; $c++) {
            $val = $row[$c]
            if ($null -ne $val) {
                $sheet.Cells.Item($r + 1, $c + 1).Value = $val
            }
        }
    }
}

Fill-Wide $xlwings
Fill-Wide $converted

# give the first new sheet the printer-settings-backed page setup
$xlwings.PageSetup.PaperSize = 9
$xlwings.PageSetup.Orientation = 1

# the trailing "转换后" sheet keeps the older-style (inch-based) margins
$converted.PageSetup.LeftMargin = 54
$converted.PageSetup.RightMargin = 54
$converted.PageSetup.TopMargin = 72
$converted.PageSetup.BottomMargin = 72
$converted.PageSetup.HeaderMargin = 36
$converted.PageSetup.FooterMargin = 36

# tidy up the original sheet's view (no longer the active tab) and keep its
# page setup as before
$orig.PageSetup.PaperSize = 9
$orig.PageSetup.Orientation = 1
[void]$orig.Range("K11").Select()

# make the new first sheet active / selected
[void]$xlwings.Activate()
[void]$xlwings.Range("A1").Select()
